$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3767.3235
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 3851.182
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 11553.546
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -11889.546

# ALC row 52
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 981.5599999999999
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 981.5599999999999
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 2944.68
$ws.Range("N52").Value = -3264.68

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4435.8335
$ws.Range("I74").Value = 4499.875
$ws.Range("J74").Value = 4307.75
$ws.Range("K74").Value = 4499.875
$ws.Range("L74").Value = 4307.75
$ws.Range("M74").Value = -3563.875
$ws.Range("N74").Value = -6179.75

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4435.8335
$ws.Range("I77").Value = 4499.875
$ws.Range("J77").Value = 4307.75
$ws.Range("K77").Value = 22499.375
$ws.Range("L77").Value = 21538.75
$ws.Range("M77").Value = -17819.375
$ws.Range("N77").Value = -30898.75

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 167.85715
$ws.Range("I96").Value = 170.83333
$ws.Range("J96").Value = 150
$ws.Range("K96").Value = 512.49999
$ws.Range("L96").Value = 450
$ws.Range("M96").Value = 860.50001
$ws.Range("N96").Value = -3196

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1857.3125
$ws.Range("I116").Value = 1679.5
$ws.Range("J116").Value = 2248.5
$ws.Range("K116").Value = 1679.5
$ws.Range("L116").Value = 2248.5
$ws.Range("M116").Value = 1762.5
$ws.Range("N116").Value = -9132.5

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1551.1091
$ws.Range("I137").Value = 1217.659
$ws.Range("J137").Value = 2884.9092
$ws.Range("K137").Value = 3652.977
$ws.Range("L137").Value = 8654.7276
$ws.Range("M137").Value = -1102.977
$ws.Range("N137").Value = -13754.7276

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10214.973
$ws.Range("I61").Value = 6334.6665
$ws.Range("J61").Value = 15307.875
$ws.Range("K61").Value = 6334.6665
$ws.Range("L61").Value = 15307.875
$ws.Range("M61").Value = -6122.6665
$ws.Range("N61").Value = -15731.875

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3877.282
$ws.Range("I74").Value = 1329.3513
$ws.Range("J74").Value = 51014
$ws.Range("K74").Value = 1329.3513
$ws.Range("L74").Value = 51014
$ws.Range("M74").Value = -455.3513
$ws.Range("N74").Value = -52762

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3877.282
$ws.Range("I77").Value = 1329.3513
$ws.Range("J77").Value = 51014
$ws.Range("K77").Value = 6646.7565
$ws.Range("L77").Value = 255070
$ws.Range("M77").Value = -2278.7565
$ws.Range("N77").Value = -263806

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2066.6667
$ws.Range("I88").Value = 2100
$ws.Range("J88").Value = 2000
$ws.Range("K88").Value = 2100
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = -1694
$ws.Range("N88").Value = -2812

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2066.6667
$ws.Range("I91").Value = 2100
$ws.Range("J91").Value = 2000
$ws.Range("K91").Value = 2100
$ws.Range("L91").Value = 2000
$ws.Range("M91").Value = -696
$ws.Range("N91").Value = -4808

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1766762.9
$ws.Range("I102").Value = 2648137.2
$ws.Range("J102").Value = 4014.2856
$ws.Range("K102").Value = 2648137.2
$ws.Range("L102").Value = 4014.2856
$ws.Range("M102").Value = -2646515.2
$ws.Range("N102").Value = -7258.2856

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2397.5
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2936
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 8808
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -13708

# ARM row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 31000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 31000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 31000
$ws.Range("N135").Value = -41140

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 10214.973
$ws.Range("I136").Value = 6334.6665
$ws.Range("J136").Value = 15307.875
$ws.Range("K136").Value = 19003.9995
$ws.Range("L136").Value = 45923.625
$ws.Range("M136").Value = -16453.9995
$ws.Range("N136").Value = -51023.625

# BSM row 55
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 60300
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 60300
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 60300
$ws.Range("N55").Value = -60846

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1126.5
$ws.Range("I94").Value = 1126
$ws.Range("J94").Value = 1127.5
$ws.Range("K94").Value = 1126
$ws.Range("L94").Value = 1127.5
$ws.Range("M94").Value = -675
$ws.Range("N94").Value = -2029.5

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1110.5
$ws.Range("I107").Value = 912.5
$ws.Range("J107").Value = 1506.5
$ws.Range("K107").Value = 912.5
$ws.Range("L107").Value = 1506.5
$ws.Range("M107").Value = 1007.5
$ws.Range("N107").Value = -5346.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 42728.56
$ws.Range("I134").Value = 2791.125
$ws.Range("J134").Value = 113728.445
$ws.Range("K134").Value = 8373.375
$ws.Range("L134").Value = 341185.335
$ws.Range("M134").Value = -5838.375
$ws.Range("N134").Value = -346255.335

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2203.8462
$ws.Range("I31").Value = 1669.8
$ws.Range("J31").Value = 3303.353
$ws.Range("K31").Value = 1669.8
$ws.Range("L31").Value = 3303.353
$ws.Range("M31").Value = -1374.8
$ws.Range("N31").Value = -3893.353

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2203.8462
$ws.Range("I34").Value = 1669.8
$ws.Range("J34").Value = 3303.353
$ws.Range("K34").Value = 1669.8
$ws.Range("L34").Value = 3303.353
$ws.Range("M34").Value = -1467.8
$ws.Range("N34").Value = -3707.353

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1492072.9
$ws.Range("I58").Value = 2115307
$ws.Range("J58").Value = 3235.9443
$ws.Range("K58").Value = 2115307
$ws.Range("L58").Value = 3235.9443
$ws.Range("M58").Value = -2115104
$ws.Range("N58").Value = -3641.9443

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4658.7144
$ws.Range("I99").Value = 4153
$ws.Range("J99").Value = 5333
$ws.Range("K99").Value = 4153
$ws.Range("L99").Value = 5333
$ws.Range("M99").Value = -2655
$ws.Range("N99").Value = -8329

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 7841.2
$ws.Range("I122").Value = 8176.5
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 24529.5
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -22079.5
$ws.Range("N122").Value = -24400

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4658.7144
$ws.Range("I126").Value = 4153
$ws.Range("J126").Value = 5333
$ws.Range("K126").Value = 12459
$ws.Range("L126").Value = 15999
$ws.Range("M126").Value = -9989
$ws.Range("N126").Value = -20939

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1492072.9
$ws.Range("I136").Value = 2115307
$ws.Range("J136").Value = 3235.9443
$ws.Range("K136").Value = 6345921
$ws.Range("L136").Value = 9707.832900000001
$ws.Range("M136").Value = -6343371
$ws.Range("N136").Value = -14807.8329

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1532
$ws.Range("I117").Value = 29
$ws.Range("J117").Value = 2033
$ws.Range("K117").Value = 87
$ws.Range("L117").Value = 6099
$ws.Range("M117").Value = 3355
$ws.Range("N117").Value = -12983

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1602997.1
$ws.Range("I139").Value = 2819310.8
$ws.Range("J139").Value = 2584.3157
$ws.Range("K139").Value = 8457932.399999999
$ws.Range("L139").Value = 7752.9471
$ws.Range("M139").Value = -8452792.399999999
$ws.Range("N139").Value = -18032.9471

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1856.9231
$ws.Range("I97").Value = 2034.2858
$ws.Range("J97").Value = 1650
$ws.Range("K97").Value = 2034.2858
$ws.Range("L97").Value = 1650
$ws.Range("M97").Value = -1538.2858
$ws.Range("N97").Value = -2642

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1028.8
$ws.Range("I93").Value = 839.2
$ws.Range("J93").Value = 1408
$ws.Range("K93").Value = 839.2
$ws.Range("L93").Value = 1408
$ws.Range("M93").Value = 408.8
$ws.Range("N93").Value = -3904

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6691.25
$ws.Range("I122").Value = 6245.2383
$ws.Range("J122").Value = 7731.9443
$ws.Range("K122").Value = 18735.7149
$ws.Range("L122").Value = 23195.8329
$ws.Range("M122").Value = -16285.7149
$ws.Range("N122").Value = -28095.8329

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3704.3809
$ws.Range("I136").Value = 2007.5238
$ws.Range("J136").Value = 7098.095
$ws.Range("K136").Value = 6022.5714
$ws.Range("L136").Value = 21294.285
$ws.Range("M136").Value = -3472.5714
$ws.Range("N136").Value = -26394.285
